$d = $word.ActiveDocument
$d.Content.Find.Execute("Windows Phone App", $true, $false, $false, $false, $false, $true, 1, $false, "Medex", 2)
